# Insert two new data rows (190 and 191) into the Albahaca price table,
# pushing the existing rows 190..266 down to 192..268.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("190:191").Insert()

# --- New row 190 ---
$ws.Cells.Item(190, 1).Value2 = 6
$ws.Cells.Item(190, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(190, 3).Value2 = "Metropolitana"
$ws.Cells.Item(190, 4).Value2 = 44510
$ws.Cells.Item(190, 5).Value2 = 13
$ws.Cells.Item(190, 6).Value2 = 100112052
$ws.Cells.Item(190, 7).Value2 = "Albahaca"
$ws.Cells.Item(190, 8).Value2 = "Sin especificar"
$ws.Cells.Item(190, 9).Value2 = "Primera"
$ws.Cells.Item(190, 10).Value2 = 500
$ws.Cells.Item(190, 11).Value2 = 5000
$ws.Cells.Item(190, 12).Value2 = 6000
$ws.Cells.Item(190, 13).Value2 = 5400
$ws.Cells.Item(190, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(190, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(190, 16).Value2 = 900
$ws.Cells.Item(190, 17).Value2 = 6
$ws.Cells.Item(190, 18).Value2 = "Hortaliza"

# --- New row 191 ---
$ws.Cells.Item(191, 1).Value2 = 6
$ws.Cells.Item(191, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(191, 3).Value2 = "Metropolitana"
$ws.Cells.Item(191, 4).Value2 = 44510
$ws.Cells.Item(191, 5).Value2 = 13
$ws.Cells.Item(191, 6).Value2 = 100112052
$ws.Cells.Item(191, 7).Value2 = "Albahaca"
$ws.Cells.Item(191, 8).Value2 = "Sin especificar"
$ws.Cells.Item(191, 9).Value2 = "Primera"
$ws.Cells.Item(191, 10).Value2 = 110
$ws.Cells.Item(191, 11).Value2 = 4000
$ws.Cells.Item(191, 12).Value2 = 5000
$ws.Cells.Item(191, 13).Value2 = 4455
$ws.Cells.Item(191, 14).Value2 = "`$/paquete"
$ws.Cells.Item(191, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(191, 16).Value2 = 4455
$ws.Cells.Item(191, 17).Value2 = 1
$ws.Cells.Item(191, 18).Value2 = "Hortaliza"

# Make sure the new date cells use the same date/time number format as the
# rest of column D (style index 2 in the original workbook).
$ws.Cells.Item(190, 4).NumberFormat = $ws.Cells.Item(192, 4).NumberFormat
$ws.Cells.Item(191, 4).NumberFormat = $ws.Cells.Item(192, 4).NumberFormat
